$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("G1").Value = "Waiter ID"
$ws.Range("H1").Value = "Chef ID"
$ws.Range("I1").Value = "Driver ID"

# Fill in zeros for the new columns on the existing rows (2-4)
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0

# Add new row 5 for a new "takeout" order (waiter order)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "takeout"
$ws.Range("C5").Value = "[1, 1, 1, 2, 2, 2]"
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = "InProgress"
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0

$null = $ws.Range("L8").Select()
